$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 970
$ws.Range("J121").Value = 970
$ws.Range("L121").Value = 2910
$ws.Range("N121").Value = -6404
$ws.Range("H132").Value = 315714.12
$ws.Range("I132").Value = 348270.94
$ws.Range("J132").Value = 998.3333
$ws.Range("K132").Value = 1044812.82
$ws.Range("L132").Value = 2994.9999
$ws.Range("M132").Value = -1042282.82
$ws.Range("N132").Value = -8054.9999
$ws.Range("H137").Value = 25001884
$ws.Range("I137").Value = 1261.1562
$ws.Range("J137").Value = 125004376
$ws.Range("K137").Value = 3783.4686
$ws.Range("L137").Value = 375013128
$ws.Range("M137").Value = -1233.4686
$ws.Range("N137").Value = -375018228

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1153
$ws.Range("I45").Value = 1037.3334
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1037.3334
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -660.3334
$ws.Range("N45").Value = -2254

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2186.1875
$ws.Range("I20").Value = 2455.4443
$ws.Range("J20").Value = 1840
$ws.Range("K20").Value = 2455.4443
$ws.Range("L20").Value = 1840
$ws.Range("M20").Value = -2208.4443
$ws.Range("N20").Value = -2334

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1262.8096
$ws.Range("I31").Value = 1262.8096
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1262.8096
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -967.8096
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1262.8096
$ws.Range("I34").Value = 1262.8096
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1262.8096
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1060.8096
$ws.Range("N34").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 804.4375
$ws.Range("I14").Value = 804.4375
$ws.Range("K14").Value = 2413.3125
$ws.Range("M14").Value = -2240.3125
$ws.Range("H68").Value = 955.2857
$ws.Range("I68").Value = 706.8333
$ws.Range("J68").Value = 1141.625
$ws.Range("K68").Value = 2120.4999
$ws.Range("L68").Value = 3424.875
$ws.Range("M68").Value = -1309.4999
$ws.Range("N68").Value = -5046.875
$ws.Range("H71").Value = 955.2857
$ws.Range("I71").Value = 706.8333
$ws.Range("J71").Value = 1141.625
$ws.Range("K71").Value = 6361.4997
$ws.Range("L71").Value = 10274.625
$ws.Range("M71").Value = -2305.4997
$ws.Range("N71").Value = -18386.625
$ws.Range("H80").Value = 8048418
$ws.Range("I80").Value = 18107616
$ws.Range("J80").Value = 1058.6
$ws.Range("K80").Value = 54322848
$ws.Range("L80").Value = 3175.8
$ws.Range("M80").Value = -54321912
$ws.Range("N80").Value = -5047.799999999999
$ws.Range("H83").Value = 8048418
$ws.Range("I83").Value = 18107616
$ws.Range("J83").Value = 1058.6
$ws.Range("K83").Value = 162968544
$ws.Range("L83").Value = 9527.4
$ws.Range("M83").Value = -162963864
$ws.Range("N83").Value = -18887.4
$ws.Range("H86").Value = 1090.4
$ws.Range("J86").Value = 1156
$ws.Range("L86").Value = 3468
$ws.Range("N86").Value = -5840
$ws.Range("H89").Value = 1090.4
$ws.Range("J89").Value = 1156
$ws.Range("L89").Value = 10404
$ws.Range("N89").Value = -22260
$ws.Range("H92").Value = 1547.8
$ws.Range("I92").Value = 525.4286
$ws.Range("J92").Value = 3933.3333
$ws.Range("K92").Value = 1576.2858
$ws.Range("L92").Value = 11799.9999
$ws.Range("M92").Value = -328.2857999999999
$ws.Range("N92").Value = -14295.9999
$ws.Range("H94").Value = 1274
$ws.Range("I94").Value = 1274
$ws.Range("K94").Value = 3822
$ws.Range("M94").Value = -3146
$ws.Range("H96").Value = 141411410
$ws.Range("J96").Value = 141411410
$ws.Range("L96").Value = 424234230
$ws.Range("N96").Value = -424238348
$ws.Range("H98").Value = 2000
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H105").Value = 908000000
$ws.Range("J105").Value = 908000000
$ws.Range("L105").Value = 2724000000
$ws.Range("N105").Value = -2724005242
$ws.Range("H107").Value = 62909.062
$ws.Range("J107").Value = 42116.125
$ws.Range("L107").Value = 126348.375
$ws.Range("N107").Value = -130188.375
$ws.Range("H110").Value = 3514.2856
$ws.Range("J110").Value = 3660
$ws.Range("L110").Value = 10980
$ws.Range("N110").Value = -19160
$ws.Range("H117").Value = 607.25
$ws.Range("I117").Value = 607.25
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1821.75
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1620.25
$ws.Range("N117").ClearContents()
$ws.Range("H129").Value = 929.3158
$ws.Range("J129").Value = 1095.5
$ws.Range("L129").Value = 3286.5
$ws.Range("N129").Value = -13286.5
$ws.Range("H131").Value = 3535.7896
$ws.Range("J131").Value = 2445.8823
$ws.Range("L131").Value = 7337.646900000001
$ws.Range("N131").Value = -17417.6469

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2733.9333
$ws.Range("I122").Value = 1913.0714
$ws.Range("J122").Value = 3452.1875
$ws.Range("K122").Value = 5739.2142
$ws.Range("L122").Value = 10356.5625
$ws.Range("M122").Value = -3289.2142
$ws.Range("N122").Value = -15256.5625

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1633.909
$ws.Range("I93").Value = 1612.1666
$ws.Range("K93").Value = 1612.1666
$ws.Range("H136").Value = 1186.619
$ws.Range("I136").Value = 961.8
$ws.Range("J136").Value = 1748.6666
$ws.Range("K136").Value = 2885.4
$ws.Range("L136").Value = 5245.9998
$ws.Range("M136").Value = -335.3999999999996
$ws.Range("N136").Value = -10345.9998

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1689.871
$ws.Range("I122").Value = 1246.3846
$ws.Range("J122").Value = 3996
$ws.Range("K122").Value = 3739.1538
$ws.Range("L122").Value = 11988
$ws.Range("M122").Value = -1289.1538
$ws.Range("N122").Value = -16888
$ws.Range("H136").Value = 10504.137
$ws.Range("I136").Value = 12399.223
$ws.Range("J136").Value = 1976.25
$ws.Range("K136").Value = 37197.669
$ws.Range("L136").Value = 5928.75
$ws.Range("M136").Value = -34647.669
$ws.Range("N136").Value = -11028.75
